{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" / \"Impact\" bullet list so each\n// line reads as an impact-focused accomplishment statement instead of a job\n// duty, and drop the old FEC-analysis / expert-testimony bullets (the list\n// shrinks from six bullets to five).\n//\n// Several of the old bullet strings also appear verbatim earlier in the\n// document (under \"PROFESSIONAL EXPERIENCE\" -> \"Partner - Siege Analytics\"),\n// so the rewrite is scoped to paragraphs that fall strictly between the\n// \"KEY ACHIEVEMENTS AND IMPACT\" Heading2 and the next Heading2\n// (\"TECHNICAL SKILLS\") to avoid touching that unrelated section.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading and the following\n// Heading 2 (the start of the next section) to bound the scan.\nlet sectionStart = -1;\nlet sectionEnd = items.length;\nfor (let i = 0; i < items.length; i++) {\n  if (sectionStart === -1) {\n    if (items[i].style === \"Heading 2\" && items[i].text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\") {\n      sectionStart = i;\n    }\n    continue;\n  }\n  if (items[i].style === \"Heading 2\") {\n    sectionEnd = i;\n    break;\n  }\n}\n\n// Old bullet text -> new bullet text. Mapping taken straight from the diff:\n// five of the six original bullets are reworded in place, and the sixth\n// (\"Provided expert testimony...\") together with its predecessor (\"Built\n// real-time FEC analysis systems...\") collapse down to one replacement\n// bullet, net-deleting a paragraph.\nconst replacements = [\n  [\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n    \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\"\n  ],\n  [\n    \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    \"\u2022 178% accuracy improvement in racial classification algorithms\"\n  ],\n  [\n    \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n    \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\"\n  ],\n  [\n    \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\",\n    \"\u2022 $4.7M savings enabled nonprofit access\"\n  ],\n  [\n    \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n    \"\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\"\n  ]\n];\nconst removedText =\n  \"\u2022 Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy\";\n\nlet removedParagraph = null;\nif (sectionStart !== -1) {\n  for (let i = sectionStart + 1; i < sectionEnd; i++) {\n    const para = items[i];\n    const text = para.text;\n    const hit = replacements.find(([oldText]) => text === oldText);\n    if (hit) {\n      para.insertText(hit[1], \"Replace\");\n      continue;\n    }\n    if (text === removedText) {\n      removedParagraph = para;\n    }\n  }\n}\n\n// Delete the now-redundant sixth bullet paragraph entirely.\nif (removedParagraph) {\n  removedParagraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" / \"Impact\" bullet list so each\n# line reads as an impact-focused accomplishment statement instead of a job\n# duty, and drop the old FEC-analysis / expert-testimony bullets (the list\n# shrinks from six bullets to five).\n#\n# Several of the old bullet strings also appear verbatim earlier in the\n# document (under \"PROFESSIONAL EXPERIENCE\" -> \"Partner - Siege Analytics\"),\n# so the rewrite is scoped to paragraphs that fall strictly between the\n# \"KEY ACHIEVEMENTS AND IMPACT\" Heading 2 and the next Heading 2\n# (\"TECHNICAL SKILLS\") to avoid touching that unrelated section.\n\n$d = $word.ActiveDocument\n\n$bullet = [char]0x2022\n\n# Old bullet text -> new bullet text, straight from the diff. Five of the six\n# original bullets are reworded in place; the sixth (\"Provided expert\n# testimony...\") together with its predecessor (\"Built real-time FEC\n# analysis systems...\") collapse down to one replacement bullet, net-deleting\n# a paragraph.\n$replacements = @(\n    @(\n        $bullet + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n        $bullet + \" Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\"\n    ),\n    @(\n        $bullet + \" Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\",\n        $bullet + \" 178% accuracy improvement in racial classification algorithms\"\n    ),\n    @(\n        $bullet + \" Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n        $bullet + \" Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\"\n    ),\n    @(\n        $bullet + \" Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" + [char]0x00B1 + \"4.2% to \" + [char]0x00B1 + \"2.1%\",\n        $bullet + \" `$4.7M savings enabled nonprofit access\"\n    ),\n    @(\n        $bullet + \" Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\",\n        $bullet + \" Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\"\n    )\n)\n$removedText = $bullet + \" Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy\"\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading and the following\n# Heading 2 (the start of the next section) to bound the scan.\n$sectionStart = -1\n$sectionEnd = $d.Paragraphs.Count + 1\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    $idx = $idx + 1\n    $styleName = $p.Style.NameLocal\n    if ($sectionStart -eq -1) {\n        $headingText = $p.Range.Text.TrimEnd()\n        if (($styleName -eq \"Heading 2\") -and ($headingText -eq \"KEY ACHIEVEMENTS AND IMPACT\")) {\n            $sectionStart = $idx\n        }\n        continue\n    }\n    if ($styleName -eq \"Heading 2\") {\n        $sectionEnd = $idx\n        break\n    }\n}\n\n$toDelete = $null\nif ($sectionStart -ne -1) {\n    $idx = 0\n    foreach ($p in $d.Paragraphs) {\n        $idx = $idx + 1\n        if ($idx -le $sectionStart) { continue }\n        if ($idx -ge $sectionEnd) { break }\n\n        $text = $p.Range.Text.TrimEnd()\n        $matched = $false\n        foreach ($pair in $replacements) {\n            if ($text -eq $pair[0]) {\n                $p.Range.Text = $pair[1]\n                $matched = $true\n                break\n            }\n        }\n        if ($matched) { continue }\n        if ($text -eq $removedText) {\n            $toDelete = $idx\n        }\n    }\n}\n\n# Delete the now-redundant sixth bullet paragraph entirely (do this last so\n# the index captured above is still valid).\nif ($toDelete -ne $null) {\n    $p = $d.Paragraphs.Item($toDelete)\n    $p.Range.Delete()\n}\n"}
